$wb = $excel.ActiveWorkbook

# This script updates static market-price snapshot values (columns H:N) across
# several "Leve" profit-tracking worksheets, per the scheduled market-data refresh.

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 756.8182
$ws.Range("I5").Value = 474.7143
$ws.Range("J5").Value = 1250.5
$ws.Range("K5").Value = 474.7143
$ws.Range("L5").Value = 1250.5
$ws.Range("M5").Value = -359.7143
$ws.Range("N5").Value = -1480.5
$ws.Range("H17").Value = 886.7432
$ws.Range("J17").Value = 871.49316
$ws.Range("L17").Value = 2614.47948
$ws.Range("N17").Value = -2950.47948
$ws.Range("H43").Value = 1947
$ws.Range("I43").Value = 1940
$ws.Range("J43").Value = 1948.75
$ws.Range("K43").Value = 1940
$ws.Range("L43").Value = 1948.75
$ws.Range("M43").Value = -1871
$ws.Range("N43").Value = -2086.75
$ws.Range("H64").Value = 4239.3
$ws.Range("I64").Value = 3599.75
$ws.Range("K64").Value = 3599.75
$ws.Range("M64").Value = -3351.75
$ws.Range("H67").Value = 4239.3
$ws.Range("I67").Value = 3599.75
$ws.Range("K67").Value = 3599.75
$ws.Range("M67").Value = -2741.75
$ws.Range("H76").Value = 181427.86
$ws.Range("H79").Value = 181427.86
$ws.Range("H88").Value = 15449.2
$ws.Range("J88").Value = 16388.555
$ws.Range("L88").Value = 16388.555
$ws.Range("N88").Value = -17200.555
$ws.Range("H91").Value = 15449.2
$ws.Range("J91").Value = 16388.555
$ws.Range("L91").Value = 16388.555
$ws.Range("N91").Value = -19196.555
$ws.Range("H101").Value = 522
$ws.Range("J101").Value = 1490
$ws.Range("L101").Value = 4470
$ws.Range("N101").Value = -7714
$ws.Range("H132").Value = 9395.75
$ws.Range("I132").Value = 5957.926
$ws.Range("K132").Value = 17873.778
$ws.Range("M132").Value = -15343.778
$ws.Range("H138").Value = 20409728
$ws.Range("I138").Value = 25642304
$ws.Range("J138").Value = 2678.8
$ws.Range("K138").Value = 76926912
$ws.Range("L138").Value = 8036.400000000001
$ws.Range("M138").Value = -76921772
$ws.Range("N138").Value = -18316.4

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H40").Value = 20000
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H45").Value = 8958.190000000001
$ws.Range("I45").Value = 13560.909
$ws.Range("J45").Value = 3895.2
$ws.Range("K45").Value = 13560.909
$ws.Range("L45").Value = 3895.2
$ws.Range("M45").Value = -13183.909
$ws.Range("N45").Value = -4649.2
$ws.Range("H97").Value = 41535.31
$ws.Range("I97").Value = 2283.2856
$ws.Range("J97").Value = 206393.8
$ws.Range("K97").Value = 2283.2856
$ws.Range("L97").Value = 206393.8
$ws.Range("M97").Value = -1787.2856
$ws.Range("N97").Value = -207385.8
$ws.Range("H132").Value = 2995.111
$ws.Range("I132").Value = 2573.5293
$ws.Range("J132").Value = 4298.1816
$ws.Range("K132").Value = 7720.5879
$ws.Range("L132").Value = 12894.5448
$ws.Range("M132").Value = -5190.5879
$ws.Range("N132").Value = -17954.5448

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2087.8918
$ws.Range("I31").Value = 1747.5454
$ws.Range("J31").Value = 2587.0667
$ws.Range("K31").Value = 1747.5454
$ws.Range("L31").Value = 2587.0667
$ws.Range("M31").Value = -1452.5454
$ws.Range("N31").Value = -3177.0667
$ws.Range("H34").Value = 2087.8918
$ws.Range("I34").Value = 1747.5454
$ws.Range("J34").Value = 2587.0667
$ws.Range("K34").Value = 1747.5454
$ws.Range("L34").Value = 2587.0667
$ws.Range("M34").Value = -1545.5454
$ws.Range("N34").Value = -2991.0667
$ws.Range("H58").Value = 6073.5415
$ws.Range("I58").Value = 6023.3823
$ws.Range("J58").Value = 6195.357
$ws.Range("K58").Value = 6023.3823
$ws.Range("L58").Value = 6195.357
$ws.Range("M58").Value = -5820.3823
$ws.Range("N58").Value = -6601.357
$ws.Range("H62").Value = 83336540
$ws.Range("I62").Value = 166669300
$ws.Range("J62").Value = 3781.3333
$ws.Range("K62").Value = 166669300
$ws.Range("L62").Value = 3781.3333
$ws.Range("M62").Value = -166668676
$ws.Range("N62").Value = -5029.3333
$ws.Range("H65").Value = 83336540
$ws.Range("I65").Value = 166669300
$ws.Range("J65").Value = 3781.3333
$ws.Range("K65").Value = 833346500
$ws.Range("L65").Value = 18906.6665
$ws.Range("M65").Value = -833343380
$ws.Range("N65").Value = -25146.6665
$ws.Range("H107").Value = 5555
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("H132").Value = 7100.2964
$ws.Range("I132").Value = 3047.4285
$ws.Range("K132").Value = 9142.2855
$ws.Range("M132").Value = -6612.2855
$ws.Range("H134").Value = 6897.4517
$ws.Range("I134").Value = 7507.9546
$ws.Range("J134").Value = 5405.1113
$ws.Range("K134").Value = 22523.8638
$ws.Range("L134").Value = 16215.3339
$ws.Range("M134").Value = -19988.8638
$ws.Range("N134").Value = -21285.3339
$ws.Range("H136").Value = 6073.5415
$ws.Range("I136").Value = 6023.3823
$ws.Range("J136").Value = 6195.357
$ws.Range("K136").Value = 18070.1469
$ws.Range("L136").Value = 18586.071
$ws.Range("M136").Value = -15520.1469
$ws.Range("N136").Value = -23686.071

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 64.82758
$ws.Range("J2").Value = 95.4375
$ws.Range("L2").Value = 572.625
$ws.Range("N2").Value = -798.625
$ws.Range("H4").Value = 3215.4
$ws.Range("I4").Value = 185
$ws.Range("K4").Value = 555
$ws.Range("M4").Value = -443
$ws.Range("H74").Value = 14999.5
$ws.Range("J74").Value = 14999.5
$ws.Range("L74").Value = 44998.5
$ws.Range("N74").Value = -47120.5
$ws.Range("H77").Value = 14999.5
$ws.Range("J77").Value = 14999.5
$ws.Range("L77").Value = 134995.5
$ws.Range("N77").Value = -145603.5
$ws.Range("H115").Value = 1890
$ws.Range("I115").Value = 1890
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 5670
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -4495
$ws.Range("N115").ClearContents()

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 21908118
$ws.Range("I80").Value = 28753244
$ws.Range("J80").Value = 3716.2
$ws.Range("K80").Value = 28753244
$ws.Range("L80").Value = 3716.2
$ws.Range("M80").Value = -28752246
$ws.Range("N80").Value = -5712.2
$ws.Range("H83").Value = 21908118
$ws.Range("I83").Value = 28753244
$ws.Range("J83").Value = 3716.2
$ws.Range("K83").Value = 143766220
$ws.Range("L83").Value = 18581
$ws.Range("M83").Value = -143761228
$ws.Range("N83").Value = -28565
$ws.Range("H126").Value = 6397.706
$ws.Range("I126").Value = 4144.6665
$ws.Range("K126").Value = 12433.9995
$ws.Range("M126").Value = -9963.999500000002
$ws.Range("H132").Value = 6948.684
$ws.Range("I132").Value = 7199.467
$ws.Range("J132").Value = 6008.25
$ws.Range("K132").Value = 21598.401
$ws.Range("L132").Value = 18024.75
$ws.Range("M132").Value = -19068.401
$ws.Range("N132").Value = -23084.75

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 719
$ws.Range("I22").Value = 705.6667
$ws.Range("J22").Value = 799
$ws.Range("K22").Value = 705.6667
$ws.Range("L22").Value = 799
$ws.Range("M22").Value = -410.6667
$ws.Range("N22").Value = -1389
$ws.Range("H27").Value = 719
$ws.Range("I27").Value = 705.6667
$ws.Range("J27").Value = 799
$ws.Range("K27").Value = 705.6667
$ws.Range("L27").Value = 799
$ws.Range("M27").Value = -598.6667
$ws.Range("N27").Value = -1013
$ws.Range("H68").Value = 3432.0715
$ws.Range("J68").Value = 15999.75
$ws.Range("L68").Value = 15999.75
$ws.Range("N68").Value = -17497.75
$ws.Range("H71").Value = 3432.0715
$ws.Range("J71").Value = 15999.75
$ws.Range("L71").Value = 79998.75
$ws.Range("N71").Value = -87486.75
$ws.Range("H93").Value = 881.875
$ws.Range("I93").Value = 982.5
$ws.Range("K93").Value = 982.5
$ws.Range("M93").Value = 265.5
$ws.Range("H122").Value = 7033.1665
$ws.Range("I122").Value = 7033.1665
$ws.Range("K122").Value = 21099.4995
$ws.Range("M122").Value = -18649.4995
$ws.Range("H132").Value = 32701.865
$ws.Range("I132").Value = 39043.2
$ws.Range("K132").Value = 117129.6
$ws.Range("M132").Value = -114599.6
$ws.Range("H140").Value = 99243.75
$ws.Range("J140").Value = 99243.75
$ws.Range("L140").Value = 99243.75
$ws.Range("N140").Value = -109603.75

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 8074.091
$ws.Range("J4").Value = 10816
$ws.Range("L4").Value = 10816
$ws.Range("N4").Value = -11042
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("H62").Value = 10561.571
$ws.Range("J62").Value = 10457.75
$ws.Range("L62").Value = 10457.75
$ws.Range("N62").Value = -11705.75
$ws.Range("H65").Value = 10561.571
$ws.Range("J65").Value = 10457.75
$ws.Range("L65").Value = 52288.75
$ws.Range("N65").Value = -58528.75
$ws.Range("H122").Value = 13605.29
$ws.Range("I122").Value = 12396.385
$ws.Range("K122").Value = 37189.155
$ws.Range("M122").Value = -34739.155
$ws.Range("H126").Value = 5690109.5
$ws.Range("I126").Value = 9619824
$ws.Range("K126").Value = 28859472
$ws.Range("M126").Value = -28857002
$ws.Range("H132").Value = 4029.9707
$ws.Range("I132").Value = 2841.8965
$ws.Range("K132").Value = 8525.6895
$ws.Range("M132").Value = -5995.6895

